$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Complexity" column (F), rows 3-25 are Fibonacci-like story
# points; the "0.5" values are stored as text. Fill F5 first so the
# "0.5" shared string is registered before "Complexity", matching the
# original authoring order. Pre-format as Text so "0.5" is kept as a
# literal value, then switch the display format to the date-ish
# "d-mmm" number format (numFmtId 16) Excel applies to this kind of
# ambiguous value - this mirrors the exact style produced upstream.
# Finally, propagate the resulting text value to F14/F24 via
# copy / paste-values so they pick up the text without inheriting that
# special number format.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0.5"
$ws.Range("F5").NumberFormat = "d-mmm"
$ws.Range("F5").Copy()
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("F24").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Header for new "Complexity" column
$ws.Range("F2").Value = "Complexity"

# Numeric complexity values for rows 3-25 (Fibonacci-like story points)
$ws.Range("F3").Value = 8
$ws.Range("F4").Value = 13
$ws.Range("F6").Value = 20
$ws.Range("F7").Value = 20
$ws.Range("F8").Value = 20
$ws.Range("F9").Value = 5
$ws.Range("F10").Value = 8
$ws.Range("F11").Value = 8
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 8
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 13
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 8
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = 5
$ws.Range("F25").Value = 1

$ws.Range("G3").Select()
